$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric (e.g. "1.000", "14.93") are
# stored as text, matching the source data which uses dotted price strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.846.90'
$ws.Range("E2").Value = '  -0.24%  '

$ws.Range("D3").Value = '1.737.87'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '240.54'
$ws.Range("E5").Value = '  +4.49%  '

$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("E7").Value = '  -0.51%  '

$ws.Range("D8").Value = '0.2749'
$ws.Range("E8").Value = '  -0.89%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.06165'
$ws.Range("E9").Value = '  +0.34%  '

$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.743.18'
$ws.Range("E10").Value = '  +0.23%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '0.07195'
$ws.Range("E11").Value = '  +1.54%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '14.93'
$ws.Range("E12").Value = '  -1.49%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.6422'
$ws.Range("E13").Value = '  -0.78%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '4.624'
$ws.Range("E14").Value = '  +2.08%  '

$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '77.51'
$ws.Range("E15").Value = '  +0.55%  '

$ws.Range("B16").Value = 'Dai'
$ws.Range("C16").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = '0.9995'
$ws.Range("E17").Value = '  +0.09%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '25.875.35'
$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").Value = '11.70'
$ws.Range("E19").Value = '  +1.41%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000006765'
$ws.Range("E20").Value = '  +1.51%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '1.968.49'
$ws.Range("E21").Value = '  +0.51%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '4.281'
$ws.Range("E22").Value = '  +0.70%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '8.642'
$ws.Range("E23").Value = '  -1.93%  '

$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '5.283'
$ws.Range("E24").Value = '  +2.31%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '138.38'
$ws.Range("E25").Value = '  -1.34%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '1.515'
$ws.Range("E26").Value = '  +0.24%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '15.20'
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '1.767'
$ws.Range("E28").Value = '  -1.26%  '

$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '105.80'
$ws.Range("E29").Value = '  +3.45%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '3.916'
$ws.Range("E30").Value = '  +5.17%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.08291'
$ws.Range("E31").Value = '  -0.48%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '3.680'
$ws.Range("E32").Value = '  +3.95%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.04632'
$ws.Range("E33").Value = '  +2.88%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.641'
$ws.Range("E34").Value = '  +1.22%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '0.9875'
$ws.Range("E35").Value = '  +1.52%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6194'
$ws.Range("E36").Value = '  -0.16%  '

$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '2.680'
$ws.Range("E37").Value = '  -0.12%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01606'
$ws.Range("E38").Value = '  +1.92%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '1.940'
$ws.Range("E39").Value = '  +1.69%  '

$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '0.9992'
$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = '97.97'
$ws.Range("E41").Value = '  -2.13%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.3846'
$ws.Range("E42").Value = '  -0.23%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.7393'
$ws.Range("E43").Value = '  +1.86%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.001'
$ws.Range("E44").Value = '  -0.47%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '0.1131'
$ws.Range("E45").Value = '  +0.89%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '6.236'
$ws.Range("E46").Value = '  +0.28%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05243'
$ws.Range("E47").Value = '  -1.57%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '54.77'
$ws.Range("E48").Value = '  +2.28%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '30.49'
$ws.Range("E49").Value = '  +1.41%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '7.569'
$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3408'
$ws.Range("E51").Value = '  +0.13%  '

# Remove the temporary text number-format so styling matches the original
# (unstyled) data cells.
$ws.Range("D2:D51").ClearFormats()
